$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Materials")

# New "Quality Planed Timber Merchants" block (rows 23-29)
$ws1.Range("B23").Value = "Quality Planed Timber Merchants in the UK | Timbersource"

$ws1.Range("B24").Value = "Ash"
$ws1.Range("C24").Value = 40

$ws1.Range("B25").Value = "Beech"
$ws1.Range("C25").Value = 42

$ws1.Range("B26").Value = "Maple"
$ws1.Range("C26").Value = 47

$ws1.Range("B27").Value = "Meranti"
$ws1.Range("C27").Value = 40

$ws1.Range("B28").Value = "American Oak"
$ws1.Range("C28").Value = 56

$ws1.Range("B29").Value = "Tulipwood"
$ws1.Range("C29").Value = 30

# Prices use the same currency format as the Metal cost column above
$ws1.Range("C4").Copy()
$ws1.Range("C24:C29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Hyperlink on the header cell, pointing at the Timbersource catalogue page
$ws1.Hyperlinks.Add($ws1.Range("B23"), "https://www.timbersource.co.uk/planed-timber", [Type]::Missing, [Type]::Missing, "Quality Planed Timber Merchants in the UK | Timbersource")
$ws1.Range("B23").Style = "Hyperlink"

# Make Materials the active sheet again, selecting the cell below the new table
$ws1.Activate() | Out-Null
$ws1.Range("B30").Select() | Out-Null
